$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in row 1
# (Set C1 before B1 so the new shared-string table entries are appended in
# the same order seen in the target workbook: "Retorno aulas presenciais"
# then "Reabertura econômica".)
$ws.Range("C1").Value = "Retorno aulas presenciais"
$ws.Range("B1").Value = "Reabertura econômica"

# Move the active selection (cosmetic change seen in the diff)
$ws.Range("J12").Select()
